$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("ECs", "Adm", "Ramp2", "ECs", 2, 0.6666666666666666, 10.793823, 32.381469, 0.2977942722284143, 0.2977942722284143, 2, 0.6666666666666666, 53.87858066666666, 161.635742, 0.5277355398832225, 0.5277355398832226, 581.5558632072218, 5234.002768864997, 0.1571566210285935, 0.1571566210285936),
    @("ECs", "Adm", "Ramp2", "FAPs", 2, 0.6666666666666666, 10.793823, 32.381469, 0.2977942722284143, 0.2977942722284143, 3, 1, 39.10044333333334, 117.30133, 0.3829850994006638, 0.3829850994006639, 422.04326456153, 3798.38938105377, 0.1140507689503476, 0.1140507689503476),
    @("ECs", "Adm", "Ramp2", "sCs", 2, 0.6666666666666666, 10.793823, 32.381469, 0.2977942722284143, 0.2977942722284143, 3, 1, 9.114878333333333, 27.344635, 0.08927936071611353, 0.08927936071611356, 98.38438339653499, 885.4594505688149, 0.02658688224947311, 0.02658688224947311),
    @("FAPs", "Adm", "Ramp2", "ECs", 3, 1, 23.80412266666667, 71.412368, 0.6567396357672255, 0.6567396357672255, 2, 0.6666666666666666, 53.87858066666666, 161.635742, 0.5277355398832225, 0.5277355398832226, 1282.532343295228, 11542.79108965706, 0.3465848462443277, 0.3465848462443277),
    @("FAPs", "Adm", "Ramp2", "FAPs", 3, 1, 23.80412266666667, 71.412368, 0.6567396357672255, 0.6567396357672255, 3, 1, 39.10044333333334, 117.30133, 0.3829850994006638, 0.3829850994006639, 930.7517494277157, 8376.765744849441, 0.2515214946846666, 0.2515214946846667),
    @("FAPs", "Adm", "Ramp2", "sCs", 3, 1, 23.80412266666667, 71.412368, 0.6567396357672255, 0.6567396357672255, 3, 1, 9.114878333333333, 27.344635, 0.08927936071611353, 0.08927936071611356, 216.9716819384089, 1952.74513744568, 0.05863329483823114, 0.05863329483823116),
    @("sCs", "Adm", "Ramp2", "ECs", 3, 1, 1.647959666666667, 4.943879000000001, 0.0454660920043603, 0.0454660920043603, 2, 0.6666666666666666, 53.87858066666666, 161.635742, 0.5277355398832225, 0.5277355398832226, 88.78972783591313, 799.1075505232182, 0.02399407261030135, 0.02399407261030136),
    @("sCs", "Adm", "Ramp2", "FAPs", 3, 1, 1.647959666666667, 4.943879000000001, 0.0454660920043603, 0.0454660920043603, 3, 1, 39.10044333333334, 117.30133, 0.3829850994006638, 0.3829850994006639, 64.43595356211891, 579.9235820590701, 0.01741283576564966, 0.01741283576564966),
    @("sCs", "Adm", "Ramp2", "sCs", 3, 1, 1.647959666666667, 4.943879000000001, 0.0454660920043603, 0.0454660920043603, 3, 1, 9.114878333333333, 27.344635, 0.08927936071611353, 0.08927936071611356, 15.02095185990722, 135.188566739165, 0.004059183628409288, 0.00405918362840929)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
